$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the trigger list in column C (rows 8-21): "13" -> "11" inside the JSON-like array string.
$newValue = '["01", "04", "08", "11", "17", "22", "26", "30"]'
$ws.Range("C8:C21").Value = $newValue

# Reflect the view state captured in the saved file: C9 selected (the cell the user edited).
$ws.Range("C9").Select()
